$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2, shifting old rows 2-7 down to 3-8.
$ws.Rows("2").Insert()
# The inserted row inherited row 1's bold/bordered formatting; clear it so it
# matches the plain (unstyled) formatting used by the other data rows.
$ws.Range("A2:H2").ClearFormats()

# Move the old header text (previously row 1) down into the new row 2.
# Leave E2, G2 and H2 untouched (they stay as blank cells, same as the diff).
$ws.Range("A2").Value = "Numberof Pieces"
$ws.Range("B2").Value = "Screw Sizes Included"
$ws.Range("C2").Value = "SpecificationsMet"
$ws.Range("D2").Value = "ContainerType"
$ws.Range("F2").Value = "Each"

# Replace row 1 with the new numeric header row (0-7), keeping its style.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
